$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Donchian")

# Rename the "index" column header (A1) to "i"
$ws.Cells.Item(1, 1).Value2 = "i"

# Convert the index column from 1-based to 0-based (decrement every data row by 1)
$lastRow = 503
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 - 1
}

# Adjust column A width to fit new narrower header/content ("i" vs "index")
$ws.Columns.Item(1).ColumnWidth = 3.14
